$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.752.16"
$ws.Range("E2").Value = "  -5.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.268.89"
$ws.Range("E3").Value = "  -6.58%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.46"
$ws.Range("E5").Value = "  -4.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.78"
$ws.Range("E6").Value = "  -4.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E8").Value = "  -4.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.263.38"
$ws.Range("E9").Value = "  -6.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.183"
$ws.Range("E10").Value = "  -10.25%  "
$ws.Range("E11").Value = "  -6.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.16"
$ws.Range("E12").Value = "  -8.77%  "
$ws.Range("E13").Value = "  -7.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.60"
$ws.Range("E14").Value = "  -6.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "631.79"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.790.07"
$ws.Range("E16").Value = "  -6.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.04"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.651.83"
$ws.Range("E18").Value = "  -5.71%  "
$ws.Range("E19").Value = "  -3.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.272.75"
$ws.Range("E20").Value = "  -6.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.30"
$ws.Range("E21").Value = "  -8.39%  "
$ws.Range("E22").Value = "  -5.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.38"
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "106.64"
$ws.Range("E24").Value = "  +7.75%  "
$ws.Range("E25").Value = "  -6.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.96"
$ws.Range("E26").Value = "  -7.53%  "
$ws.Range("E27").Value = "  -7.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.51"
$ws.Range("E28").Value = "  -5.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.67"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.23"
$ws.Range("E30").Value = "  -7.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.93"
$ws.Range("E31").Value = "  -7.41%  "
$ws.Range("E32").Value = "  -7.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.02"
$ws.Range("E33").Value = "  -5.41%  "
$ws.Range("E34").Value = "  -4.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.64"
$ws.Range("E35").Value = "  -6.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.729.44"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "521.17"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("E39").Value = "  -4.63%  "
$ws.Range("E40").Value = "  -7.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.131"
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("E42").Value = "  -7.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.40"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "32.90"
$ws.Range("E44").Value = "  -4.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.336"
$ws.Range("E45").Value = "  -10.09%  "
$ws.Range("E46").Value = "  -3.55%  "
$ws.Range("E47").Value = "  -6.81%  "
$ws.Range("E48").Value = "  -4.05%  "
$ws.Range("E49").Value = "  -8.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("E51").Value = "  +0.58%  "
